# "Store module has been added" — add a new "Sheet2" (Store/Order data)
# after the existing "Sheet1" (Userdata), matching the Petstore "Store"
# schema: userID, petID, Quantity, shipDate, status, complete.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet immediately after Sheet1 so it lands at the end
# and naturally becomes the active/selected tab (matching activeTab=1
# and sheet1 losing tabSelected).
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"

# ---- Header row -------------------------------------------------------
# Build the font/style table in the same order the original author's
# Excel session produced it: plain 14pt first (header row), then bold
# 14pt (most body cells), then the two vertically-centered variants.
$ws2.Range("A1").Value = "userID"
$ws2.Range("B1").Value = "petID"
$ws2.Range("C1").Value = "Quantity"
$ws2.Range("E1").Value = "status"
$ws2.Range("F1").Value = "complete"
$ws2.Range("D1").Value = "shipDate"
$ws2.Range("A1:F1").Font.Size = 14

# ---- Body rows ----------------------------------------------------------
$ws2.Range("A2").Value = 5
$ws2.Range("B2").Value = 10
$ws2.Range("C2").Value = 122
$ws2.Range("D2").Value = "2023-12-06T04:03:05.234Z"
$ws2.Range("E2").Value = "placed"
$ws2.Range("F2").Value = $true

$ws2.Range("A3").Value = 6
$ws2.Range("B3").Value = 11
$ws2.Range("C3").Value = 234
$ws2.Range("D3").Value = "2023-12-06T04:03:05.234Z"
$ws2.Range("E3").Value = "placed"
$ws2.Range("F3").Value = $true

$ws2.Range("A4").Value = 7
$ws2.Range("B4").Value = 12
$ws2.Range("C4").Value = 222
$ws2.Range("D4").Value = "2023-12-06T04:03:05.234Z"
$ws2.Range("E4").Value = "placed"
$ws2.Range("F4").Value = $true

# Most of the data cells: bold 14pt (reuses the plain-14pt font created
# above as an intermediate step, then derives the bold variant). Touch
# one such range first so the bold font/style are allocated before the
# vertically-centered variants below (keeps font/xf allocation order in
# the workbook identical to the original authoring session).
$ws2.Range("B2:C3").Font.Size = 14
$ws2.Range("B2:C3").Font.Bold = $true
$ws2.Range("E2:F3").Font.Size = 14
$ws2.Range("E2:F3").Font.Bold = $true
$ws2.Range("A4:C4").Font.Size = 14
$ws2.Range("A4:C4").Font.Bold = $true
$ws2.Range("E4:F4").Font.Size = 14
$ws2.Range("E4:F4").Font.Bold = $true

# A2:A3 — bold 14pt, vertically centered.
$ws2.Range("A2:A3").Font.Size = 14
$ws2.Range("A2:A3").Font.Bold = $true
$ws2.Range("A2:A3").VerticalAlignment = -4108

# D2:D4 (shipDate) — 14pt, light-green text, vertically centered. Done
# last so its new font/style land after the bold-centered ones above.
$ws2.Range("D2:D4").Font.Size = 14
$ws2.Range("D2:D4").Font.Color = 10681506
$ws2.Range("D2:D4").VerticalAlignment = -4108

# ---- Layout -------------------------------------------------------------
$ws2.Rows("1:4").RowHeight = 18.75

$ws2.Range("A1").ColumnWidth = 16.584
$ws2.Range("B1").ColumnWidth = 18.751
$ws2.Range("C1").ColumnWidth = 16.417
$ws2.Range("D1").ColumnWidth = 33.25
$ws2.Range("E1").ColumnWidth = 14.584
$ws2.Range("F1").ColumnWidth = 16.417

$ws2.PageSetup.Orientation = 1

$ws2.Range("D2").Select()
